$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  (columns: A=title/meta B=Hora_Scrap C=Hora_Llegada D=Linea
#                    E=Minutos F=Parada G=Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 08:56:26"
$ws1.Range("A3").Value = "Total filas: 95"

$rows1 = @(
    @("09:02", "215A_EL PATO", 6),
    @("09:04", "11_ETCHEVERRY", 8),
    @("09:08", "23_HERNANDEZ", 12),
    @("09:11", "16_P MOR-SANTA ANA", 15),
    @("09:13", "10_OLMOS", 17),
    @("09:17", "27_EL RETIRO", 21),
    @("09:21", "26_HERNANDEZ", 25),
    @("09:23", "16_SANTA ANA", 27),
    @("09:24", "11_ETCHEVERRY", 28),
    @("09:32", "15_ABASTO", 36),
    @("09:33", "10_OLMOS", 37),
    @("09:35", "23_HERNANDEZ", 39),
    @("09:42", "215C_EL PATO", 46),
    @("09:44", "14_ABASTO", 48),
    @("09:52", "15_ABASTO", 56),
    @("09:53", "10_OLMOS", 57),
    @("10:04", "11_ETCHEVERRY", 68),
    @("10:11", "16_P MOR-SANTA ANA", 75),
    @("10:12", "15_ABASTO", 76),
    @("10:22", "17_ROMERO", 86),
    @("10:27", "215A_EL PATO", 91)
)

$r = 76
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 2).Value = "08:56:16"
    $ws1.Cells.Item($r, 3).Value = $row[0]
    $ws1.Cells.Item($r, 4).Value = $row[1]
    $ws1.Cells.Item($r, 5).Value = $row[2]
    $ws1.Cells.Item($r, 6).Value = "LP1912"
    $ws1.Cells.Item($r, 7).Value = "30/12/2025"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215  (columns: A=title/meta B=Fecha C=Hora_Scrap
#                        D=Hora_Llegada E=Linea F=Minutos G=Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 30/12/2025 08:56:26"
$ws2.Range("A3").Value = "Total filas: 13"

$rows2 = @(
    @("09:02", "215A_EL PATO", 6),
    @("09:42", "215C_EL PATO", 46),
    @("10:27", "215A_EL PATO", 91)
)

$r = 12
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 2).Value = "30/12/2025"
    $ws2.Cells.Item($r, 3).Value = "08:56:16"
    $ws2.Cells.Item($r, 4).Value = $row[0]
    $ws2.Cells.Item($r, 5).Value = $row[1]
    $ws2.Cells.Item($r, 6).Value = $row[2]
    $ws2.Cells.Item($r, 7).Value = "LP1912"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173  (columns: A=title/meta B=Fecha C=Hora_Scrap
#                        D=Hora_Llegada E=Linea F=Minutos G=Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 08:56:26"
$ws3.Range("A3").Value = "Total filas: 13"

$ws3.Cells.Item(13, 2).Value = "30/12/2025"
$ws3.Cells.Item(13, 3).Value = "08:56:21"
$ws3.Cells.Item(13, 4).Value = "09:09"
$ws3.Cells.Item(13, 5).Value = "215D_LA PLATA"
$ws3.Cells.Item(13, 6).Value = 13
$ws3.Cells.Item(13, 7).Value = "L6203"

$ws3.Cells.Item(14, 2).Value = "30/12/2025"
$ws3.Cells.Item(14, 3).Value = "08:56:26"
$ws3.Cells.Item(14, 4).Value = "10:03"
$ws3.Cells.Item(14, 5).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(14, 6).Value = 67
$ws3.Cells.Item(14, 7).Value = "L6173"
